$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "generated on" timestamp in the header banner
$ws.Range("A1").Value = "Reporte generado el 19/05/2025 a las 00:16"

# Two brand-new inventory rows are appended below the existing data (rows 12 and 13).
# Create them first by copying the formatting (style s="3") of an existing data row
# so the new cells inherit the same borders/alignment, then fill in their values below.
$ws.Range("A10:F10").Copy()
$ws.Range("A12:F13").PasteSpecial(-4122)

# Row 3 ("HPe") was removed from the inventory; every row below it (old rows 4-11)
# moves up one position for columns B:F, while column A (the sequential ID) is left as-is.
# Row 11 and the two newly appended rows (12, 13) get their own fresh product data.

$ws.Cells.Item(3, 2).Value = "Cama Queen"
$ws.Cells.Item(3, 3).Value = 10
$ws.Cells.Item(3, 4).Value = 699
$ws.Cells.Item(3, 5).Value = "unidad"
$ws.Cells.Item(3, 6).Value = "Dormitorio"

$ws.Cells.Item(4, 2).Value = "Juego de terraza"
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 899.9
$ws.Cells.Item(4, 5).Value = "juego"
$ws.Cells.Item(4, 6).Value = "Exteriores"

$ws.Cells.Item(5, 2).Value = "Lámpara de pie"
$ws.Cells.Item(5, 3).Value = 15
$ws.Cells.Item(5, 4).Value = 85.75
$ws.Cells.Item(5, 5).Value = "unidad"
$ws.Cells.Item(5, 6).Value = "Sala"

$ws.Cells.Item(6, 2).Value = "Escritorio ejecutivo"
$ws.Cells.Item(6, 3).Value = 10000000000
$ws.Cells.Item(6, 4).Value = 399
$ws.Cells.Item(6, 5).Value = "unidad"
$ws.Cells.Item(6, 6).Value = "Oficina"

$ws.Cells.Item(7, 2).Value = "Velador doble"
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(7, 4).Value = 120
$ws.Cells.Item(7, 5).Value = "par"
$ws.Cells.Item(7, 6).Value = "Dormitorio"

$ws.Cells.Item(8, 2).Value = "Silla comedor"
$ws.Cells.Item(8, 3).Value = 32000
$ws.Cells.Item(8, 4).Value = 45.999
$ws.Cells.Item(8, 5).Value = "unidad"
$ws.Cells.Item(8, 6).Value = "Exteriores"

$ws.Cells.Item(9, 2).Value = "Hola"
$ws.Cells.Item(9, 3).Value = 21
$ws.Cells.Item(9, 4).Value = 123
$ws.Cells.Item(9, 5).Value = "unidad"
$ws.Cells.Item(9, 6).Value = "Comedor"

$ws.Cells.Item(10, 2).Value = "Hola 123"
$ws.Cells.Item(10, 3).Value = 12
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = "unidad"
$ws.Cells.Item(10, 6).Value = "Sala"

$ws.Cells.Item(11, 2).Value = "Sdf213"
$ws.Cells.Item(11, 3).Value = 23
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = "juego"
$ws.Cells.Item(11, 6).Value = "Comedor"

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Sad"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 23
$ws.Cells.Item(12, 5).Value = "juego"
$ws.Cells.Item(12, 6).Value = "Oficina"

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Sda"
$ws.Cells.Item(13, 3).Value = 12
$ws.Cells.Item(13, 4).Value = 2131232
$ws.Cells.Item(13, 5).Value = "par"
$ws.Cells.Item(13, 6).Value = "Comedor"
